$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$col = $used.Columns.Item(7)   # Column G ("Recorded By")

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$cell = $col.Find($target, [Type]::Missing, [Type]::Missing, 1)
if ($cell -ne $null) {
    $firstRow = $cell.Row
    while ($cell -ne $null) {
        $cell.Value = $replacement
        $cell = $col.FindNext($cell)
        if ($cell -eq $null) { break }
        if ($cell.Row -eq $firstRow) { break }
    }
}
